$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "DONE" status text is changed to "d". Update every cell that currently
# holds "DONE" (B1:B77) at once so the shared string is edited in place
# instead of creating a duplicate string entry.
$ws.Range("B1:B77").Value = "d"

# Clear out the status column for all rows except the header/first row.
$ws.Range("B2:B77").ClearContents()

# Insert a new row under the first row, pushing all the existing data down
# by one, and duplicate the first row's A value into the newly inserted row.
$ws.Rows.Item(2).Insert()
$ws.Range("A1").Copy($ws.Range("A2"))

# Update the active selection to B2, matching the saved view state.
$ws.Range("B2").Select()
